# Updated code quality rules and rel rating
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "BannedPaths" rule row (old row 35). Everything
# below shifts up by one.
$ws.Range("A35").EntireRow.Delete()

# Insert a fresh row in its place (pushes "AEM Rules:AEM-3" back down to
# row 41) and populate it with the renamed/re-rated rule.
$ws.Range("A40").EntireRow.Insert()
$ws.Range("A40").Value = "BannedPath"
$ws.Range("B40").Value = "Customer packages should not install content under /libs"
$ws.Range("C40").Value = "Bug"
$ws.Range("D40").Value = "Critical"

# Update the remembered selection to match the author's final cursor
# position.
$ws.Range("A37").Select()
